# energy.xlsx help-text / style update
#  - updates to energy help (Loads/Protocol tab copy, "Save Defaults" rename)
#  - correct electric co2 conversion value (Value 100% pressure explanation)
#  - remove debug prints (N/A to data - handled via text only)
#
# Sheet order in the workbook: 1=Intro, 2=Details, 3=Loads, 4=Protocol

$wb = $excel.ActiveWorkbook

$wsIntro    = $wb.Worksheets.Item("Intro")
$wsDetails  = $wb.Worksheets.Item("Details")
$wsLoads    = $wb.Worksheets.Item("Loads")
$wsProtocol = $wb.Worksheets.Item("Protocol")

# A17 (Protocol): "Save as Defaults" -> "Save Defaults", highlighted in red
$wsProtocol.Range("A17").Value = "Save Defaults"
$wsProtocol.Range("A17").Font.Color = 255

# --- Loads tab -------------------------------------------------------------

# A2: intro paragraph rewritten, and no longer uses the wrapped-text style
$wsLoads.Range("A2").Value = "tn:Begin by making entries on the Loads tab  to define the sources of energy used by this roast. It might be a good idea to save those settings as defaults to be used to calculate the energy consumption of future roasts"
$wsLoads.Range("A2").WrapText = $false

# B17: "Value 100%" description rewritten (correct electric co2 conversion value)
$wsLoads.Range("B17").Value = "When an Event is selected this value can be set to match the 100% load setting to the event setting.  This is useful when the 100% load setting is recorded as a different number in the Event.  For instance, maybe the burner event is recorded as 10x the kPa reading on the gas manometer.  An event value of 35 is recoded to signify 3.5 kPa, which is 50% pressure.  If the 100% burner setting corresponds to 7 kPa then the 'Value 100%' should be set to 70, which is 7 * 10  = 70.  Thus 3.5 kPa will be seen by he energy calculator as 50%.  For pressure readings be sure to tick the Pressure box.  Heat energy readings are normally 0%-100% and do not require any adjustment to this  setting."

# A19: "Save as Defaults" -> "Save Defaults"
$wsLoads.Range("A19").Value = "Save Defaults"

# --- Protocol tab ------------------------------------------------------------

# B11: "Between Batches" row description rewritten
$wsProtocol.Range("B11").Value = "This row sets the values for between batches protocol for the roasting session.  Percentage or measured values may be entered for each burner.  When a percentage is used the Duration field must be set.\n\nBetween Batches energy is applied to each batch of the roasting session, except the first batch.  Tick the 'Between Batches after Pre-Heating' box to apply Between Batches energies to the first batch of the session too."

# --- Active tab / selection bookkeeping ------------------------------------
# Workbook now opens on the Loads tab (3rd tab, 0-indexed activeTab="2"),
# with Loads!A2 and Protocol!A17 the last-selected cells.

$wsLoads.Activate() | Out-Null
$wsLoads.Range("A2").Select() | Out-Null

$wsProtocol.Activate() | Out-Null
$wsProtocol.Range("A17").Select() | Out-Null

$wsLoads.Activate() | Out-Null
